$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4240
$ws.Range("B3").Value = 2280.1
$ws.Range("B4").Value = 681.5
$ws.Range("B5").Value = 245.6
$ws.Range("B6").Value = 1115.1
$ws.Range("B7").Value = 312.1
$ws.Range("B8").Value = 8874.400000000001
